$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Excel enum constants used below (PasteSpecial "paste" arg).
$xlPasteFormats = -4122

# --- Pre-format the text columns (B,D,E,F) for rows 2-22 as Text so that
# date-like / numeric-like strings ("2025-04-03", "000486", ...) are not
# auto-converted to dates/numbers when the .Value is assigned.
$ws.Range("B2:B22").NumberFormat = "@"
$ws.Range("D2:D22").NumberFormat = "@"
$ws.Range("E2:E22").NumberFormat = "@"
$ws.Range("F2:F22").NumberFormat = "@"

# Row 2
$ws.Cells.Item(2, 1).Value = 2
$ws.Cells.Item(2, 2).Value = '2025-03-24'
$ws.Cells.Item(2, 3).Value = 150
$ws.Cells.Item(2, 4).Value = 'JURUA ESTALEIROS E NAVEGACAO LTDA'
$ws.Cells.Item(2, 5).Value = '000088'
$ws.Cells.Item(2, 6).Value = 'VASSOURA PIACAVA 20 FUROS'
$ws.Cells.Item(2, 7).Value = 192
$ws.Cells.Item(2, 8).Value = $false

# Row 3
$ws.Cells.Item(3, 1).Value = 0
$ws.Cells.Item(3, 2).Value = '2025-03-25'
$ws.Cells.Item(3, 3).Value = 300
$ws.Cells.Item(3, 4).Value = 'MAP SERVICOS DE CONSERVACAO - EIRELI'
$ws.Cells.Item(3, 5).Value = '000098'
$ws.Cells.Item(3, 6).Value = 'PANO DE CHAO FLANELADO C REFORCADO ITATEX 42x62CM'
$ws.Cells.Item(3, 7).Value = 1865
$ws.Cells.Item(3, 8).Value = $false

# Row 4
$ws.Cells.Item(4, 1).Value = 19
$ws.Cells.Item(4, 2).Value = '2025-03-25'
$ws.Cells.Item(4, 3).Value = 200
$ws.Cells.Item(4, 4).Value = 'MAP SERVICOS DE CONSERVACAO - EIRELI'
$ws.Cells.Item(4, 5).Value = '000486'
$ws.Cells.Item(4, 6).Value = 'NAFTALINA 40G PCT RUBI'
$ws.Cells.Item(4, 7).Value = 298
$ws.Cells.Item(4, 8).Value = $false

# Row 5
$ws.Cells.Item(5, 1).Value = 9
$ws.Cells.Item(5, 2).Value = '2025-03-25'
$ws.Cells.Item(5, 3).Value = 15
$ws.Cells.Item(5, 4).Value = 'AMAZONIA REFEICOES E SERVICOS LTDA'
$ws.Cells.Item(5, 5).Value = '000924'
$ws.Cells.Item(5, 6).Value = 'COPO POTE DESCARTAVEL TRANSP 100ML CX/20'
$ws.Cells.Item(5, 7).Value = 57
$ws.Cells.Item(5, 8).Value = $false

# Row 6
$ws.Cells.Item(6, 1).Value = 5
$ws.Cells.Item(6, 2).Value = '2025-03-25'
$ws.Cells.Item(6, 3).Value = 40
$ws.Cells.Item(6, 4).Value = 'MM DA AMAZONIA INDUSTRIA E COMERCIO DE PLASTICOS LTDA.'
$ws.Cells.Item(6, 5).Value = '000288'
$ws.Cells.Item(6, 6).Value = 'TOUCA DESCARTAVEL TNT TALGE PCT C/ 100 UND'
$ws.Cells.Item(6, 7).Value = 646
$ws.Cells.Item(6, 8).Value = $false

# Row 7
$ws.Cells.Item(7, 1).Value = 10
$ws.Cells.Item(7, 2).Value = '2025-03-25'
$ws.Cells.Item(7, 3).Value = 96
$ws.Cells.Item(7, 4).Value = 'MM DA AMAZONIA INDUSTRIA E COMERCIO DE PLASTICOS LTDA.'
$ws.Cells.Item(7, 5).Value = '000032'
$ws.Cells.Item(7, 6).Value = 'LIMPADOR VEJA MULTIUSO GOLD 500ML'
$ws.Cells.Item(7, 7).Value = 1148
$ws.Cells.Item(7, 8).Value = $false

# Row 8
$ws.Cells.Item(8, 1).Value = 1
$ws.Cells.Item(8, 2).Value = '2025-03-25'
$ws.Cells.Item(8, 3).Value = 120
$ws.Cells.Item(8, 4).Value = 'MAP SERVICOS DE CONSERVACAO - EIRELI'
$ws.Cells.Item(8, 5).Value = '000350'
$ws.Cells.Item(8, 6).Value = 'DESODORISADOR LADY AEROSSOL 360 ML LAVANDA'
$ws.Cells.Item(8, 7).Value = 856
$ws.Cells.Item(8, 8).Value = $true

# Row 9
$ws.Cells.Item(9, 1).Value = 20
$ws.Cells.Item(9, 2).Value = '2025-03-25'
$ws.Cells.Item(9, 3).Value = 96
$ws.Cells.Item(9, 4).Value = 'MANJAR SERVICOS GERAIS SA'
$ws.Cells.Item(9, 5).Value = '000583'
$ws.Cells.Item(9, 6).Value = 'LIMPA ALUMINIO BRINORT 500ML'
$ws.Cells.Item(9, 7).Value = 99
$ws.Cells.Item(9, 8).Value = $false

# Row 10
$ws.Cells.Item(10, 1).Value = 7
$ws.Cells.Item(10, 2).Value = '2025-03-26'
$ws.Cells.Item(10, 3).Value = 10
$ws.Cells.Item(10, 4).Value = 'JURUA ESTALEIROS E NAVEGACAO LTDA'
$ws.Cells.Item(10, 5).Value = '000425'
$ws.Cells.Item(10, 6).Value = 'COADOR DE CAFE EG (EXTRA GRANDE)'
$ws.Cells.Item(10, 7).Value = -4
$ws.Cells.Item(10, 8).Value = $false

# Row 11
$ws.Cells.Item(11, 1).Value = 14
$ws.Cells.Item(11, 2).Value = '2025-03-26'
$ws.Cells.Item(11, 3).Value = 200
$ws.Cells.Item(11, 4).Value = 'JURUA ESTALEIROS E NAVEGACAO LTDA'
$ws.Cells.Item(11, 5).Value = '000122'
$ws.Cells.Item(11, 6).Value = 'SABAO EM PO ALA LAVANDA ROUPAS 400G'
$ws.Cells.Item(11, 7).Value = 363
$ws.Cells.Item(11, 8).Value = $false

# Row 12
$ws.Cells.Item(12, 1).Value = 4
$ws.Cells.Item(12, 2).Value = '2025-03-26'
$ws.Cells.Item(12, 3).Value = 40
$ws.Cells.Item(12, 4).Value = 'V V REFEICOES LTDA'
$ws.Cells.Item(12, 5).Value = '000091'
$ws.Cells.Item(12, 6).Value = 'VASSOURA VARRE CANTO COM CABO PLASTIFICADO'
$ws.Cells.Item(12, 7).Value = 138
$ws.Cells.Item(12, 8).Value = $false

# Row 13
$ws.Cells.Item(13, 1).Value = 18
$ws.Cells.Item(13, 2).Value = '2025-03-28'
$ws.Cells.Item(13, 3).Value = 350
$ws.Cells.Item(13, 4).Value = 'AMAZONPEL COMERCIO DE MATERIAIS DE LIMPEZA LTDA'
$ws.Cells.Item(13, 5).Value = '000494'
$ws.Cells.Item(13, 6).Value = 'FIBRA DE LIMPEZA PESADA 98X229MM SLIM NOBRE'
$ws.Cells.Item(13, 7).Value = 331
$ws.Cells.Item(13, 8).Value = $true

# Row 14
$ws.Cells.Item(14, 1).Value = 12
$ws.Cells.Item(14, 2).Value = '2025-04-01'
$ws.Cells.Item(14, 3).Value = 20
$ws.Cells.Item(14, 4).Value = 'SAT BRAS INDUSTRIA ELETRONICA DA AMAZONIA LTDA.'
$ws.Cells.Item(14, 5).Value = '000258'
$ws.Cells.Item(14, 6).Value = 'DISPENSER PAPEL HIGIENICO ROLAO 300-500M BRANCO NOBRE STREET'
$ws.Cells.Item(14, 7).Value = 40
$ws.Cells.Item(14, 8).Value = $true

# Row 15
$ws.Cells.Item(15, 1).Value = 17
$ws.Cells.Item(15, 2).Value = '2025-04-01'
$ws.Cells.Item(15, 3).Value = 61
$ws.Cells.Item(15, 4).Value = 'AMAZONIA FORMULA LTDA'
$ws.Cells.Item(15, 5).Value = '000799'
$ws.Cells.Item(15, 6).Value = 'SACO DE LIXO 30L REFORCADO PACOTINHO C/10 UND'
$ws.Cells.Item(15, 7).Value = 361
$ws.Cells.Item(15, 8).Value = $false

# Row 16
$ws.Cells.Item(16, 1).Value = 6
$ws.Cells.Item(16, 2).Value = '2025-04-02'
$ws.Cells.Item(16, 3).Value = 250
$ws.Cells.Item(16, 4).Value = 'RH MULTI SERVICOS ADMINISTRATIVOS S.A'
$ws.Cells.Item(16, 5).Value = '000041'
$ws.Cells.Item(16, 6).Value = 'LUVAS DESCARTAVEIS C/ 100 UND'
$ws.Cells.Item(16, 7).Value = 1637
$ws.Cells.Item(16, 8).Value = $false

# Row 17
$ws.Cells.Item(17, 1).Value = 8
$ws.Cells.Item(17, 2).Value = '2025-04-02'
$ws.Cells.Item(17, 3).Value = 100
$ws.Cells.Item(17, 4).Value = 'REVEMAR COMERCIO DE MOTOS LTDA'
$ws.Cells.Item(17, 5).Value = '010041'
$ws.Cells.Item(17, 6).Value = 'PAPEL A4 REPORT RESMA C 500FLS'
$ws.Cells.Item(17, 7).Value = 46
$ws.Cells.Item(17, 8).Value = $false

# Row 18
$ws.Cells.Item(18, 1).Value = 3
$ws.Cells.Item(18, 2).Value = '2025-04-02'
$ws.Cells.Item(18, 3).Value = 50
$ws.Cells.Item(18, 4).Value = 'REVEMAR COMERCIO DE MOTOS LTDA'
$ws.Cells.Item(18, 5).Value = '000054'
$ws.Cells.Item(18, 6).Value = 'PAPEL HIGIENICO 8X300 NEWPAPER 100% Celulose'
$ws.Cells.Item(18, 7).Value = -41
$ws.Cells.Item(18, 8).Value = $false

# Row 19
$ws.Cells.Item(19, 1).Value = 11
$ws.Cells.Item(19, 2).Value = '2025-04-02'
$ws.Cells.Item(19, 3).Value = 30
$ws.Cells.Item(19, 4).Value = 'REVEMAR COMERCIO DE MOTOS LTDA'
$ws.Cells.Item(19, 5).Value = '000434'
$ws.Cells.Item(19, 6).Value = 'FRASCO COM VALVULA PUMP 450ML (p/alcool/sabonete) - NOBRE'
$ws.Cells.Item(19, 7).Value = 35
$ws.Cells.Item(19, 8).Value = $true

# Row 20
$ws.Cells.Item(20, 1).Value = 13
$ws.Cells.Item(20, 2).Value = '2025-04-02'
$ws.Cells.Item(20, 3).Value = 150
$ws.Cells.Item(20, 4).Value = 'REVEMAR COMERCIO DE MOTOS LTDA'
$ws.Cells.Item(20, 5).Value = '000349'
$ws.Cells.Item(20, 6).Value = 'DESODORISADOR LADY AEROSSOL 360ML TALCO SUAVE CARINHO'
$ws.Cells.Item(20, 7).Value = 876
$ws.Cells.Item(20, 8).Value = $true

# Row 21
$ws.Cells.Item(21, 1).Value = 15
$ws.Cells.Item(21, 2).Value = '2025-04-02'
$ws.Cells.Item(21, 3).Value = 100
$ws.Cells.Item(21, 4).Value = 'REVEMAR COMERCIO DE MOTOS LTDA'
$ws.Cells.Item(21, 5).Value = '000057'
$ws.Cells.Item(21, 6).Value = 'PAPEL TOALHA 8x100M NEWPAPER 100% CELULOSE'
$ws.Cells.Item(21, 7).Value = -123
$ws.Cells.Item(21, 8).Value = $false

# Row 22
$ws.Cells.Item(22, 1).Value = 16
$ws.Cells.Item(22, 2).Value = '2025-04-03'
$ws.Cells.Item(22, 3).Value = 120
$ws.Cells.Item(22, 4).Value = 'SIND.DAS EMPRESAS DE TRANSP.DE PASSAGEIROS DO EST.'
$ws.Cells.Item(22, 5).Value = '000782'
$ws.Cells.Item(22, 6).Value = 'ALCOOL LIQUIDO 70% INPM SANTA CRUZ 1L'
$ws.Cells.Item(22, 7).Value = 72
$ws.Cells.Item(22, 8).Value = $false

# --- Restore plain (non-text-forced) formatting on the text columns so the
# cells end up with the same default styling as the rest of the sheet (no
# leftover "@" number format), by pasting formats only from a plain cell.
$ws.Range("C3").Copy() | Out-Null
$ws.Range("B2:B22").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("D2:D22").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("E2:E22").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("F2:F22").PasteSpecial($xlPasteFormats) | Out-Null

# --- New row 22, column A needs the same bordered/bold/centered style used
# by A2:A21 ("Dia" index column) - copy format from A21.
$ws.Range("A21").Copy() | Out-Null
$ws.Range("A22").PasteSpecial($xlPasteFormats) | Out-Null

# --- Update the worksheet dimension / used range to include the new row.
$ws.Application.CutCopyMode = $false

Write-Output "edit applied"